# Update timestamp column (Z) values on the active sheet to reflect the
# new run's timestamps (the notebook was re-run, producing new timestamps
# for each logged row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2  = "2025-10-17T07:09:29.810110"
    3  = "2025-10-17T07:09:29.810110"
    4  = "2025-10-17T07:09:29.810110"
    5  = "2025-10-17T07:09:29.810110"
    6  = "2025-10-17T07:09:29.810110"
    7  = "2025-10-17T07:09:29.810110"
    8  = "2025-10-17T07:09:29.810110"
    9  = "2025-10-17T07:09:29.810110"
    10 = "2025-10-17T07:09:29.810110"
    11 = "2025-10-17T07:09:29.810110"
    12 = "2025-10-17T07:09:29.810110"
    13 = "2025-10-17T07:09:29.810110"
    14 = "2025-10-17T07:09:29.810110"
    15 = "2025-10-17T07:09:29.810110"
    16 = "2025-10-17T07:09:29.888446"
    17 = "2025-10-17T07:09:29.888446"
    18 = "2025-10-17T07:09:29.889719"
    19 = "2025-10-17T07:09:29.890143"
    20 = "2025-10-17T07:09:29.890143"
    21 = "2025-10-17T07:09:29.890143"
    22 = "2025-10-17T07:09:29.890143"
    23 = "2025-10-17T07:09:29.890143"
    24 = "2025-10-17T07:09:29.890655"
    25 = "2025-10-17T07:09:29.890655"
    26 = "2025-10-17T07:09:30.068058"
    27 = "2025-10-17T07:09:30.069056"
    28 = "2025-10-17T07:09:30.069056"
    29 = "2025-10-17T07:09:30.069056"
    30 = "2025-10-17T07:09:30.069056"
    31 = "2025-10-17T07:09:30.069056"
    32 = "2025-10-17T07:09:30.069056"
    33 = "2025-10-17T07:09:30.069056"
    34 = "2025-10-17T07:09:30.069056"
    35 = "2025-10-17T07:09:30.069056"
    36 = "2025-10-17T07:09:30.069056"
    37 = "2025-10-17T07:09:30.070057"
    38 = "2025-10-17T07:09:30.070057"
    39 = "2025-10-17T07:09:30.070057"
    40 = "2025-10-17T07:09:30.070057"
    41 = "2025-10-17T07:09:30.070057"
    42 = "2025-10-17T07:09:30.070057"
    43 = "2025-10-17T07:09:30.070057"
    44 = "2025-10-17T07:09:30.070057"
    45 = "2025-10-17T07:09:30.071057"
    46 = "2025-10-17T07:09:30.071057"
    47 = "2025-10-17T07:09:30.071057"
    48 = "2025-10-17T07:09:30.071057"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
